$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 3 by copying row 2 so it inherits the same cell styles (e.g. the
# date format on column A), then overwrite with the new row's values.
$ws.Range("A2:N2").Copy($ws.Range("A3:N3"))

$ws.Range("A3").Value = 42605.885381944441
$ws.Range("B3").Value = -34
$ws.Range("C3").Value = 43
$ws.Range("D3").Value = 56
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 80
$ws.Range("G3").Value = 25349
$ws.Range("H3").Value = 6426
$ws.Range("I3").Value = 376
$ws.Range("J3").Value = 41
$ws.Range("K3").Value = 53
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 8
$ws.Range("N3").Value = "Named"
